# Append new scraped listings and refresh the "取得日時" (fetched-at) timestamp.
# Mirrors: insert 2 new rows (new postings found this run) at positions 7 and 11,
# shift the existing ones down, stamp every data row with the new run timestamp,
# and rebuild the F-column hyperlinks so they stay aligned with their rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @("2025-09-30 01:17:23", "【急募】リスト抽出ツール開発のフリーランスを探しています!", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5402362", 128, "◆ツール,開発"),
    @("2025-09-30 01:17:23", "商標登録のための依頼者と弁理士をつなぐマッチングサイト開発", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5402790", 93, "◆開発 ◇サイト"),
    @("2025-09-30 01:17:23", "日程を作成するシステムの開発", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5402412", 78, "◆開発"),
    @("2025-09-30 01:17:23", "【急募】Snowflake IntelligenceでのAgent開発者を探しています!", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5403054", 68, "◆開発"),
    @("2025-09-30 01:17:23", "【ペットのアバター化】Pawsitiveプロトタイプ開発の依頼", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5399313", 68, "◆開発"),
    @("2025-09-30 01:17:23", "【急募】LINE WORKSで定期メッセージ配信ツール作成依頼", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5403166", 65, "◆ツール"),
    @("2025-09-30 01:17:23", "【急募】新しい口コミサイトの構築をお手伝いください!", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5402277", 38, "◇サイト"),
    @("2025-09-30 01:17:23", "【急募】ブランドサイトの新商品更新アラート作成依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5402794", 33, "◇サイト"),
    @("2025-09-30 01:17:23", "MYSQLからGoogleスプレッドシートへデータ取り込み及びスプレッドシート改修", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400606", 30, "◇MySQL"),
    @("2025-09-30 01:17:23", "初回 サブスクペイからCSVデータをダウンロードし、データベース同期するプログラムの作成", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5403072", 18, $null),
    @("2025-09-30 01:17:23", "限定公開 PR 限定公開の仕事", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5399347", 13, $null),
    @("2025-09-30 01:17:23", "エンジニア面談をお願い致します", "システム開発", "~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5402603", 10, $null)
)

# Make room for the two new postings discovered in this run (rows 7 and 11
# in the final layout) by inserting blank rows and letting the rows below
# shift down, same as Excel's own "Insert Sheet Rows" command.
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(11).Insert()

# Drop every existing hyperlink now (they get rebuilt below against the
# shifted/refreshed rows) - on this sheet a single Hyperlinks.Delete() call
# clears the whole collection.
$ws.Cells.Item(2, 6).Hyperlinks.Delete()

# Write every data row (2-13) fresh: new run timestamp in col A, the scraped
# fields in B-G, and the optional skill-summary in H (cleared entirely for
# rows that don't have one, matching the source data).
$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $col = $c + 1
        $val = $row[$c]
        if ($col -eq 8 -and $val -eq $null) {
            $ws.Cells.Item($r, $col).Value = ""
        } else {
            $ws.Cells.Item($r, $col).Value = $val
        }
    }
    $r++
}

# Re-create the F-column hyperlinks in row order (F2..F13) pointing at each
# row's URL (read back via .Value2 - .Value's getter is unreliable in this
# host), then restore the plain "Hyperlink" cell style (Hyperlinks.Add mints
# a fresh style entry otherwise).
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $target = $cell.Value2
    $ws.Hyperlinks.Add($cell, $target) | Out-Null
    $cell.Style = "Hyperlink"
}

Write-Output "applied scrape-append edit"
